$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.800.82'
$ws.Range("E2").Value = '  -3.07%  '
$ws.Range("D3").Value = '1.616.93'
$ws.Range("E3").Value = '  -3.42%  '
$ws.Range("D5").Value = '308.12'
$ws.Range("E5").Value = '  -1.96%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '0.3933'
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("D8").Value = '0.3833'
$ws.Range("E8").Value = '  -2.53%  '
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").Value = '49.34'
$ws.Range("E10").Value = '  -2.11%  '
$ws.Range("D11").Value = '1.353'
$ws.Range("E11").Value = '  -2.93%  '
$ws.Range("D12").Value = '0.08432'
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("D13").Value = '23.64'
$ws.Range("E13").Value = '  -6.73%  '
$ws.Range("D14").Value = '7.028'
$ws.Range("E14").Value = '  -3.58%  '
$ws.Range("D15").Value = '7.571'
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("D16").Value = '0.00001277'
$ws.Range("E16").Value = '  -2.76%  '
$ws.Range("D17").Value = '1.618.94'
$ws.Range("E17").Value = '  -3.25%  '
$ws.Range("D18").Value = '93.82'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = '0.06929'
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("D20").Value = '19.94'
$ws.Range("E20").Value = '  -5.91%  '
$ws.Range("D21").Value = '6.804'
$ws.Range("E21").Value = '  -3.50%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.36%  '
$ws.Range("D24").Value = '23.800.42'
$ws.Range("E24").Value = '  -3.09%  '
$ws.Range("D25").Value = '2.446'
$ws.Range("E25").Value = '  +4.26%  '
$ws.Range("D26").Value = '2.832'
$ws.Range("E26").Value = '  +2.40%  '
$ws.Range("D27").Value = '22.19'
$ws.Range("E27").Value = '  -3.48%  '
$ws.Range("D28").Value = '156.98'
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("D29").Value = '139.83'
$ws.Range("E29").Value = '  -3.81%  '
$ws.Range("D30").Value = '5.291'
$ws.Range("E30").Value = '  -9.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.800'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.84%  '
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("D33").Value = '1.793.59'
$ws.Range("E33").Value = '  -3.49%  '
$ws.Range("D34").Value = '0.08091'
$ws.Range("E34").Value = '  -1.94%  '
$ws.Range("D35").Value = '0.9755'
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("D36").Value = '0.02876'
$ws.Range("E36").Value = '  -6.44%  '
$ws.Range("E37").Value = '  -5.07%  '
$ws.Range("D38").Value = '0.2659'
$ws.Range("E38").Value = '  -4.88%  '
$ws.Range("D39").Value = '0.09124'
$ws.Range("E39").Value = '  -5.25%  '
$ws.Range("D40").Value = '10.31'
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").Value = '13.54'
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").Value = '1.424'
$ws.Range("E42").Value = '  -6.02%  '
$ws.Range("D43").Value = '0.7486'
$ws.Range("E43").Value = '  -4.66%  '
$ws.Range("D44").Value = '16.03'
$ws.Range("E44").Value = '  -2.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6900'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("D46").Value = '2.467'
$ws.Range("E46").Value = '  -3.41%  '
$ws.Range("E47").Value = '  -2.46%  '
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("D49").Value = '0.08222'
$ws.Range("D50").Value = '134.53'
$ws.Range("E50").Value = '  -2.22%  '
$ws.Range("D51").Value = '1.202'
$ws.Range("E51").Value = '  -9.21%  '
